$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.075.15'
$ws.Range("E2").Value = '  -0.36%  '
$ws.Range("D3").Value = '1.917.73'
$ws.Range("E3").Value = '  +0.73%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.58'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.81%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.0000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5045'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4025'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.40%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08253'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.109'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.04'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.06'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.52%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.432'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.15%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.912.51'
$ws.Range("E14").Value = '  +0.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.292'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.86%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9994'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.22%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.24'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.68%  '
$ws.Range("E18").Value = '  -1.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06508'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.10'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.43%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.951'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.75%  '
$ws.Range("D23").Value = '30.112.54'
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.31'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.48%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.200'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '22.47'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.81%  '
$ws.Range("D27").Value = '2.136.06'
$ws.Range("E27").Value = '  +0.46%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '161.80'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.12%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.281'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.25%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '129.14'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.128'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.84%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1038'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.007'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.87%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.784'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02447'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.90%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.355'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06428'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.14%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2163'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.71%  '
$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6532'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.44%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.773'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.26%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.201'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.98%  '
$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.44'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.99%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.222'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.90%  '
$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.200'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.76%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.30'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6039'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.28%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.640'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '123.67'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.215'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.32%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '78.68'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.61%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.128'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.96%  '
